# Laboratorio 7 - Entrega final
# Update the memory value recorded for the "Carga de Catálogo CHAINING" table
# (row 10, column C on sheet "Datos Lab7") and leave the sheet in the same
# view state (zoom + selected cell) that the author ended up with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")
$ws.Activate()

# Corrected memory reading for the n=2 CHAINING catalog-load measurement.
$ws.Range("C10").Value = 29649.014999999999

# Zoom the sheet in and leave the selection where the author left it.
$excel.ActiveWindow.Zoom = 170
$ws.Range("E15").Select() | Out-Null
